$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 0")

# Update Gasoline_Price (B), Diesel_Price (C), LPG_Price (D) values for rows 2-23 and 25-27.
# Row 24 keeps its original values (180,31 / 197,03 / 104,72) - no change there.

$ws.Range("B2").Value = "187,70 "
$ws.Range("C2").Value = "193,31 "
$ws.Range("D2").Value = "132,85 "
$ws.Range("B3").Value = "144,06 "
$ws.Range("C3").Value = "145,25 "
$ws.Range("D3").Value = "80,10 "
$ws.Range("B4").Value = "152,40 "
$ws.Range("C4").Value = "153,60 "
$ws.Range("D4").Value = "77,10 "
$ws.Range("B5").Value = "215,89 "
$ws.Range("C5").Value = "233,30 "
$ws.Range("D5").Value = "118,76 "
$ws.Range("B6").Value = "171,75 "
$ws.Range("C6").Value = "166,93 "
$ws.Range("D6").Value = "78,64 "
$ws.Range("B7").Value = "206,18 "
$ws.Range("C7").Value = "191,79 "
$ws.Range("D7").Value = "124,66 "
$ws.Range("B8").Value = "190,15 "
$ws.Range("C8").Value = "180,45 "
$ws.Range("D8").Value = "129,34 "
$ws.Range("B9").Value = "185,59 "
$ws.Range("C9").Value = "176,47 "
$ws.Range("D9").Value = "112,97 "
$ws.Range("B10").Value = "209,92 "
$ws.Range("C10").Value = "201,85 "
$ws.Range("D10").Value = "114,14 "
$ws.Range("B11").Value = "236,82 "
$ws.Range("C11").Value = "227,22 "
$ws.Range("D11").Value = "132,52 "
$ws.Range("B12").Value = "209,68 "
$ws.Range("C12").Value = "181,38 "
$ws.Range("D12").Value = "110,63 "
$ws.Range("B13").Value = "173,08 "
$ws.Range("C13").Value = "163,72 "
$ws.Range("D13").Value = "93,56 "
$ws.Range("B14").Value = "175,05 "
$ws.Range("C14").Value = "175,28 "
$ws.Range("D14").Value = "104,51 "
$ws.Range("B15").Value = "205,36 "
$ws.Range("C15").Value = "198,34 "
$ws.Range("D15").Value = "112,27 "
$ws.Range("B16").Value = "216,47 "
$ws.Range("C16").Value = "205,01 "
$ws.Range("D16").Value = "85,60 "
$ws.Range("B17").Value = "165,13 "
$ws.Range("C17").Value = "178,93 "
$ws.Range("D17").Value = "86,54 "
$ws.Range("B18").Value = "185,48 "
$ws.Range("C18").Value = "179,40 "
$ws.Range("D18").Value = "96,48 "
$ws.Range("B19").Value = "145,39 "
$ws.Range("C19").Value = "134,94 "
$ws.Range("D19").Value = "87,42 "
$ws.Range("B20").Value = "241,73 "
$ws.Range("C20").Value = "218,81 "
$ws.Range("D20").Value = "102,09 "
$ws.Range("B21").Value = "168,68 "
$ws.Range("C21").Value = "170,32 "
$ws.Range("D21").Value = "110,08 "
$ws.Range("B22").Value = "210,85 "
$ws.Range("C22").Value = "197,17 "
$ws.Range("D22").Value = "113,79 "
$ws.Range("B23").Value = "166,84 "
$ws.Range("C23").Value = "171,07 "
$ws.Range("D23").Value = "81,31 "
$ws.Range("B25").Value = "179,16 "
$ws.Range("C25").Value = "183,49 "
$ws.Range("D25").Value = "107,82 "
$ws.Range("B26").Value = "181,85 "
$ws.Range("C26").Value = "176,82 "
$ws.Range("D26").Value = "88,41 "
$ws.Range("B27").Value = "139,84 "
$ws.Range("C27").Value = "142,73 "
$ws.Range("D27").Value = "82,91 "

# Move the active selection on the "Table 0" sheet from F3 to H4.
$ws.Range("H4").Select()
